$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-7 (columns A-I), reflecting the reordered scenarios
# and cleaned-up file paths / descriptions.
$data = @(
    @{ Row=2; A="natural"; B="sheep";                 C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='["../data/grasslands/livestock/grassland_sheep_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - Only Sheep"; H="../LEAFs/SOC/rasters"; I=100 },
    @{ Row=3; A="natural"; B="goat";                  C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='["../data/grasslands/livestock/grassland_goat_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - Only Goat"; H="../LEAFs/SOC/rasters"; I=100 },
    @{ Row=4; A="natural"; B="cattle_avg_sheep_goat";  C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='[ "../data/grasslands/livestock/grassland_sheep_annual_fym.tif", "../data/grasslands/livestock/grassland_goat_annual_fym.tif", "../data/grasslands/livestock/grassland_cattle_other_avgdw_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - All Animals, including  Cattle - Low developing world cattle productivity"; H="../LEAFs/SOC/rasters"; I=100 },
    @{ Row=5; A="natural"; B="cattle_avg";             C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='["../data/grasslands/livestock/grassland_cattle_other_avgdw_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - Only Cattle - Average developing world cattle productivity"; H="../LEAFs/SOC/rasters"; I=100 },
    @{ Row=6; A="natural"; B="cattle_hps";             C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='["../data/grasslands/livestock/grassland_cattle_other_hpsdw_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - Only Cattle - High developing world cattle productivity"; H="../LEAFs/SOC/rasters"; I=100 },
    @{ Row=7; A="natural"; B="cattle_lps";             C="../data/land_use/lu_Grassland.tif"; D=14; E="../data/soil_weather/uhth_pet_locationonly.tif"; F='["../data/grasslands/livestock/grassland_cattle_other_lpsdw_annual_fym.tif"]'; G="Natural Grassland  SOC content for 2030 - Only Cattle - Low developing world cattle productivity"; H="../LEAFs/SOC/rasters"; I=100 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}

# Update the selected cell to A2 (was C8)
$ws.Range("A2").Select()
